$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.024701
$ws.Range("H2").Value = 9.074103000000001
$ws.Range("I2").Value = 0.1596375877334842
$ws.Range("J2").Value = 0.1596375877334843
$ws.Range("M2").Value = 10.20278466666667
$ws.Range("N2").Value = 30.608354
$ws.Range("O2").Value = 0.1130484251481675
$ws.Range("P2").Value = 0.1130484251481675
$ws.Range("Q2").Value = 30.86037298405133
$ws.Range("R2").Value = 277.743356856462
$ws.Range("S2").Value = 0.01804677788772282
$ws.Range("T2").Value = 0.01804677788772282
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.024701
$ws.Range("H3").Value = 9.074103000000001
$ws.Range("I3").Value = 0.1596375877334842
$ws.Range("J3").Value = 0.1596375877334843
$ws.Range("O3").Value = 0.6307851663035086
$ws.Range("P3").Value = 0.6307851663035084
$ws.Range("Q3").Value = 172.194044095878
$ws.Range("R3").Value = 1549.746396862902
$ws.Range("S3").Value = 0.1006970223267568
$ws.Range("T3").Value = 0.1006970223267568
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.024701
$ws.Range("H4").Value = 9.074103000000001
$ws.Range("I4").Value = 0.1596375877334842
$ws.Range("J4").Value = 0.1596375877334843
$ws.Range("M4").Value = 17.26138866666667
$ws.Range("N4").Value = 51.784166
$ws.Range("O4").Value = 0.1912588443635774
$ws.Range("P4").Value = 0.1912588443635773
$ws.Range("Q4").Value = 52.21053956145533
$ws.Range("R4").Value = 469.8948560530981
$ws.Range("S4").Value = 0.03053210054689539
$ws.Range("T4").Value = 0.03053210054689539
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.024701
$ws.Range("H5").Value = 9.074103000000001
$ws.Range("I5").Value = 0.1596375877334842
$ws.Range("J5").Value = 0.1596375877334843
$ws.Range("M5").Value = 5.858002
$ws.Range("N5").Value = 17.574006
$ws.Range("O5").Value = 0.06490756418474665
$ws.Range("P5").Value = 0.06490756418474665
$ws.Range("Q5").Value = 17.718704507402
$ws.Range("R5").Value = 159.468340566618
$ws.Range("S5").Value = 0.01036168697210925
$ws.Range("T5").Value = 0.01036168697210925
$ws.Range("I6").Value = 0.6072559333217162
$ws.Range("J6").Value = 0.6072559333217163
$ws.Range("M6").Value = 10.20278466666667
$ws.Range("N6").Value = 30.608354
$ws.Range("O6").Value = 0.1130484251481675
$ws.Range("P6").Value = 0.1130484251481675
$ws.Range("Q6").Value = 117.3918051829569
$ws.Range("R6").Value = 1056.526246646612
$ws.Range("S6").Value = 0.06864932692390065
$ws.Range("T6").Value = 0.06864932692390065
$ws.Range("I7").Value = 0.6072559333217162
$ws.Range("J7").Value = 0.6072559333217163
$ws.Range("O7").Value = 0.6307851663035086
$ws.Range("P7").Value = 0.6307851663035084
$ws.Range("S7").Value = 0.3830480348891311
$ws.Range("T7").Value = 0.3830480348891311
$ws.Range("I8").Value = 0.6072559333217162
$ws.Range("J8").Value = 0.6072559333217163
$ws.Range("M8").Value = 17.26138866666667
$ws.Range("N8").Value = 51.784166
$ws.Range("O8").Value = 0.1912588443635774
$ws.Range("P8").Value = 0.1912588443635773
$ws.Range("Q8").Value = 198.6071098966609
$ws.Range("R8").Value = 1787.463989069948
$ws.Range("S8").Value = 0.116143068040037
$ws.Range("T8").Value = 0.116143068040037
$ws.Range("I9").Value = 0.6072559333217162
$ws.Range("J9").Value = 0.6072559333217163
$ws.Range("M9").Value = 5.858002
$ws.Range("N9").Value = 17.574006
$ws.Range("O9").Value = 0.06490756418474665
$ws.Range("P9").Value = 0.06490756418474665
$ws.Range("Q9").Value = 67.40134698638533
$ws.Range("R9").Value = 606.612122877468
$ws.Range("S9").Value = 0.03941550346864753
$ws.Range("T9").Value = 0.03941550346864754
$ws.Range("G10").Value = 4.368554666666666
$ws.Range("H10").Value = 13.105664
$ws.Range("I10").Value = 0.2305634602787257
$ws.Range("J10").Value = 0.2305634602787257
$ws.Range("M10").Value = 10.20278466666667
$ws.Range("N10").Value = 30.608354
$ws.Range("O10").Value = 0.1130484251481675
$ws.Range("P10").Value = 0.1130484251481675
$ws.Range("Q10").Value = 44.57142256856177
$ws.Range("R10").Value = 401.142803117056
$ws.Range("S10").Value = 0.02606483608122202
$ws.Range("T10").Value = 0.02606483608122202
$ws.Range("G11").Value = 4.368554666666666
$ws.Range("H11").Value = 13.105664
$ws.Range("I11").Value = 0.2305634602787257
$ws.Range("J11").Value = 0.2305634602787257
$ws.Range("O11").Value = 0.6307851663035086
$ws.Range("P11").Value = 0.6307851663035084
$ws.Range("Q11").Value = 248.698663076864
$ws.Range("R11").Value = 2238.287967691776
$ws.Range("S11").Value = 0.1454360106354284
$ws.Range("T11").Value = 0.1454360106354284
$ws.Range("G12").Value = 4.368554666666666
$ws.Range("H12").Value = 13.105664
$ws.Range("I12").Value = 0.2305634602787257
$ws.Range("J12").Value = 0.2305634602787257
$ws.Range("M12").Value = 17.26138866666667
$ws.Range("N12").Value = 51.784166
$ws.Range("O12").Value = 0.1912588443635774
$ws.Range("P12").Value = 0.1912588443635773
$ws.Range("Q12").Value = 75.40732001291377
$ws.Range("R12").Value = 678.6658801162239
$ws.Range("S12").Value = 0.04409730096537665
$ws.Range("T12").Value = 0.04409730096537665
$ws.Range("G13").Value = 4.368554666666666
$ws.Range("H13").Value = 13.105664
$ws.Range("I13").Value = 0.2305634602787257
$ws.Range("J13").Value = 0.2305634602787257
$ws.Range("M13").Value = 5.858002
$ws.Range("N13").Value = 17.574006
$ws.Range("O13").Value = 0.06490756418474665
$ws.Range("P13").Value = 0.06490756418474665
$ws.Range("Q13").Value = 25.59100197444267
$ws.Range("R13").Value = 230.319017769984
$ws.Range("S13").Value = 0.01496531259669867
$ws.Range("T13").Value = 0.01496531259669867
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04818333333333333
$ws.Range("H14").Value = 0.14455
$ws.Range("I14").Value = 0.002543018666073676
$ws.Range("J14").Value = 0.002543018666073677
$ws.Range("M14").Value = 10.20278466666667
$ws.Range("N14").Value = 30.608354
$ws.Range("O14").Value = 0.1130484251481675
$ws.Range("P14").Value = 0.1130484251481675
$ws.Range("Q14").Value = 0.4916041745222222
$ws.Range("R14").Value = 4.4244375707
$ws.Range("S14").Value = 0.0002874842553220228
$ws.Range("T14").Value = 0.0002874842553220228
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04818333333333333
$ws.Range("H15").Value = 0.14455
$ws.Range("I15").Value = 0.002543018666073676
$ws.Range("J15").Value = 0.002543018666073677
$ws.Range("O15").Value = 0.6307851663035086
$ws.Range("P15").Value = 0.6307851663035084
$ws.Range("Q15").Value = 2.7430423783
$ws.Range("R15").Value = 24.6873814047
$ws.Range("S15").Value = 0.001604098452192211
$ws.Range("T15").Value = 0.001604098452192211
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04818333333333333
$ws.Range("H16").Value = 0.14455
$ws.Range("I16").Value = 0.002543018666073676
$ws.Range("J16").Value = 0.002543018666073677
$ws.Range("M16").Value = 17.26138866666667
$ws.Range("N16").Value = 51.784166
$ws.Range("O16").Value = 0.1912588443635774
$ws.Range("P16").Value = 0.1912588443635773
$ws.Range("Q16").Value = 0.8317112439222222
$ws.Range("R16").Value = 7.485401195300001
$ws.Range("S16").Value = 0.0004863748112682574
$ws.Range("T16").Value = 0.0004863748112682574
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04818333333333333
$ws.Range("H17").Value = 0.14455
$ws.Range("I17").Value = 0.002543018666073676
$ws.Range("J17").Value = 0.002543018666073677
$ws.Range("M17").Value = 5.858002
$ws.Range("N17").Value = 17.574006
$ws.Range("O17").Value = 0.06490756418474665
$ws.Range("P17").Value = 0.06490756418474665
$ws.Range("Q17").Value = 0.2822580630333333
$ws.Range("R17").Value = 2.540322567300001
$ws.Range("S17").Value = 0.000165061147291186
$ws.Range("T17").Value = 0.000165061147291186
